$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.855.93'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '1.710.01'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').Value = '317.42'
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('D7').Value = '0.3943'
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('D8').Value = '0.4059'
$ws.Range('E8').Value = '  -0.80%  '
$ws.Range('D9').Value = '1.493'
$ws.Range('E9').Value = '  -0.83%  '
$ws.Range('D10').Value = '0.9994'
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('D12').Value = '0.08824'
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('D13').Value = '26.57'
$ws.Range('E13').Value = '  +9.58%  '
$ws.Range('D14').Value = '7.502'
$ws.Range('E14').Value = '  -3.60%  '
$ws.Range('D15').Value = '8.139'
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('D16').Value = '0.00001362'
$ws.Range('E16').Value = '  +2.23%  '
$ws.Range('D17').Value = '1.711.47'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').Value = '96.38'
$ws.Range('E18').Value = '  -3.40%  '
$ws.Range('D19').Value = '0.07195'
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('D20').Value = '21.24'
$ws.Range('E20').Value = '  +6.06%  '
$ws.Range('D21').Value = '7.313'
$ws.Range('E21').Value = '  +1.14%  '
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.38%  '
$ws.Range('D23').Value = '14.40'
$ws.Range('E23').Value = '  -2.17%  '
$ws.Range('D24').Value = '24.843.22'
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').Value = '3.011'
$ws.Range('E25').Value = '  -3.30%  '
$ws.Range('D26').Value = '2.340'
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('D27').Value = '23.29'
$ws.Range('E27').Value = '  +0.61%  '
$ws.Range('D28').Value = '166.54'
$ws.Range('E28').Value = '  +0.85%  '
$ws.Range('D29').Value = '6.019'
$ws.Range('E29').Value = '  +15.94%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '145.33'
$ws.Range('E30').Value = '  +4.31%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '8.479'
$ws.Range('E31').Value = '  -8.45%  '
$ws.Range('B32').Value = 'WEMIXTOKEN'
$ws.Range('C32').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D32').Value = '2.245'
$ws.Range('E32').Value = '  +14.19%  '
$ws.Range('B33').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C33').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D33').Value = '1.901.66'
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('D34').Value = '0.08837'
$ws.Range('E34').Value = '  -3.24%  '
$ws.Range('D35').Value = '0.03143'
$ws.Range('E35').Value = '  +2.22%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '1.046'
$ws.Range('E36').Value = '  -3.43%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = '7.216'
$ws.Range('E37').Value = '  -10.70%  '
$ws.Range('D38').Value = '0.2884'
$ws.Range('E38').Value = '  +2.06%  '
$ws.Range('D39').Value = '11.00'
$ws.Range('E39').Value = '  -1.56%  '
$ws.Range('D40').Value = '0.8378'
$ws.Range('E40').Value = '  +6.89%  '
$ws.Range('D41').Value = '0.09242'
$ws.Range('E41').Value = '  -0.69%  '
$ws.Range('D42').Value = '14.16'
$ws.Range('E42').Value = '  -2.89%  '
$ws.Range('D43').Value = '1.477'
$ws.Range('E43').Value = '  +0.11%  '
$ws.Range('D44').Value = '17.53'
$ws.Range('E44').Value = '  +7.33%  '
$ws.Range('D45').Value = '2.693'
$ws.Range('E45').Value = '  +1.22%  '
$ws.Range('D46').Value = '0.7414'
$ws.Range('E46').Value = '  +2.02%  '
$ws.Range('D47').Value = '4.253'
$ws.Range('E47').Value = '  +0.39%  '
$ws.Range('D48').Value = '1.400'
$ws.Range('E48').Value = '  +2.50%  '
$ws.Range('D49').Value = '1.000'
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('D50').Value = '141.36'
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('D51').Value = '0.08302'
$ws.Range('E51').Value = '  +3.01%  '
